$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: insert two new columns before column D (shift D:K -> F:M)
$ws.Range("D1:E1").EntireColumn.Insert()

# Step 2: the insert copies formatting from column C (General) into the new D:E
#         columns; copy the real number formats from column F (old column D)
#         across the used data rows so D:E match the date/number styles used
#         throughout the table.
$src = $ws.Range("F5:F102")
$dst = $ws.Range("D5:E102")
$src.Copy()
$dst.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# Step 3: populate the two new quarter columns (D = period ending 2018-12-31,
#         E = period ending 2018-09-30) with the newly reported figures.
$ws.Cells.Item(7,4).Value = 43465
$ws.Cells.Item(7,5).Value = 43373
$ws.Cells.Item(8,4).Value = 288200
$ws.Cells.Item(8,5).Value = 282200
$ws.Cells.Item(9,4).Value = 41700
$ws.Cells.Item(9,5).Value = 41900
$ws.Cells.Item(10,4).Value = 246500
$ws.Cells.Item(10,5).Value = 240300
$ws.Cells.Item(12,4).Value = 10300
$ws.Cells.Item(12,5).Value = 11000
$ws.Cells.Item(13,4).Value = 0
$ws.Cells.Item(13,5).Value = 0
$ws.Cells.Item(14,4).Value = 54600
$ws.Cells.Item(14,5).Value = 16200
$ws.Cells.Item(15,4).Value = 125500
$ws.Cells.Item(15,5).Value = 107300
$ws.Cells.Item(17,4).Value = 282900
$ws.Cells.Item(17,5).Value = 140000
$ws.Cells.Item(18,4).Value = 5300
$ws.Cells.Item(18,5).Value = 142200
$ws.Cells.Item(20,4).Value = 169200
$ws.Cells.Item(20,5).Value = -34100
$ws.Cells.Item(21,4).Value = 300000
$ws.Cells.Item(21,5).Value = 215400
$ws.Cells.Item(22,4).Value = 20100
$ws.Cells.Item(22,5).Value = 20700
$ws.Cells.Item(23,4).Value = 154400
$ws.Cells.Item(23,5).Value = 87400
$ws.Cells.Item(24,4).Value = 54600
$ws.Cells.Item(24,5).Value = 22200
$ws.Cells.Item(25,4).Value = 0
$ws.Cells.Item(25,5).Value = 0
$ws.Cells.Item(26,4).Value = 99900
$ws.Cells.Item(26,5).Value = 65200
$ws.Cells.Item(27,4).Value = 91600
$ws.Cells.Item(27,5).Value = 57600
$ws.Cells.Item(28,4).Value = 0
$ws.Cells.Item(28,5).Value = 0
$ws.Cells.Item(29,4).Value = 0
$ws.Cells.Item(29,5).Value = "NA"
$ws.Cells.Item(30,4).Value = 0
$ws.Cells.Item(30,5).Value = 0
$ws.Cells.Item(31,4).Value = 0
$ws.Cells.Item(31,5).Value = 0
$ws.Cells.Item(32,4).Value = -169200
$ws.Cells.Item(32,5).Value = 34100
$ws.Cells.Item(33,4).Value = 91600
$ws.Cells.Item(33,5).Value = 57600
$ws.Cells.Item(34,4).Value = 0
$ws.Cells.Item(34,5).Value = 0
$ws.Cells.Item(35,4).Value = 91600
$ws.Cells.Item(35,5).Value = 57600
$ws.Cells.Item(38,4).Value = 43465
$ws.Cells.Item(38,5).Value = 43373
$ws.Cells.Item(41,4).Value = 235000
$ws.Cells.Item(41,5).Value = 274100
$ws.Cells.Item(42,4).Value = 0
$ws.Cells.Item(42,5).Value = 0
$ws.Cells.Item(43,4).Value = 132900
$ws.Cells.Item(43,5).Value = 151700
$ws.Cells.Item(44,4).Value = 19900
$ws.Cells.Item(44,5).Value = 26700
$ws.Cells.Item(45,4).Value = 76800
$ws.Cells.Item(45,5).Value = 13200
$ws.Cells.Item(46,4).Value = 464600
$ws.Cells.Item(46,5).Value = 465700
$ws.Cells.Item(47,4).Value = 0
$ws.Cells.Item(47,5).Value = 0
$ws.Cells.Item(48,4).Value = 3672000
$ws.Cells.Item(48,5).Value = 3617000
$ws.Cells.Item(49,4).Value = 2200
$ws.Cells.Item(49,5).Value = 56400
$ws.Cells.Item(50,4).Value = 0
$ws.Cells.Item(50,5).Value = 0
$ws.Cells.Item(51,4).Value = 0
$ws.Cells.Item(51,5).Value = 0
$ws.Cells.Item(52,4).Value = 27300
$ws.Cells.Item(52,5).Value = 19100
$ws.Cells.Item(53,4).Value = 0
$ws.Cells.Item(53,5).Value = 0
$ws.Cells.Item(54,4).Value = 4166000
$ws.Cells.Item(54,5).Value = 4158300
$ws.Cells.Item(57,4).Value = 186200
$ws.Cells.Item(57,5).Value = 200100
$ws.Cells.Item(58,4).Value = 0
$ws.Cells.Item(58,5).Value = 0
$ws.Cells.Item(59,4).Value = 216200
$ws.Cells.Item(59,5).Value = 350500
$ws.Cells.Item(60,4).Value = 402400
$ws.Cells.Item(60,5).Value = 550700
$ws.Cells.Item(61,4).Value = 1417700
$ws.Cells.Item(61,5).Value = 1422100
$ws.Cells.Item(62,4).Value = 286900
$ws.Cells.Item(62,5).Value = 212400
$ws.Cells.Item(63,4).Value = 0
$ws.Cells.Item(63,5).Value = 0
$ws.Cells.Item(64,4).Value = 0
$ws.Cells.Item(64,5).Value = 0
$ws.Cells.Item(65,4).Value = 0
$ws.Cells.Item(65,5).Value = 0
$ws.Cells.Item(66,4).Value = 2254800
$ws.Cells.Item(66,5).Value = 2329100
$ws.Cells.Item(68,4).Value = 0
$ws.Cells.Item(68,5).Value = 0
$ws.Cells.Item(69,4).Value = 0
$ws.Cells.Item(69,5).Value = 0
$ws.Cells.Item(70,4).Value = 164400
$ws.Cells.Item(70,5).Value = 162800
$ws.Cells.Item(71,4).Value = 0
$ws.Cells.Item(71,5).Value = 0
$ws.Cells.Item(72,4).Value = -375800
$ws.Cells.Item(72,5).Value = -475600
$ws.Cells.Item(73,4).Value = 0
$ws.Cells.Item(73,5).Value = 0
$ws.Cells.Item(74,4).Value = 0
$ws.Cells.Item(74,5).Value = 0
$ws.Cells.Item(75,4).Value = 0
$ws.Cells.Item(75,5).Value = 0
$ws.Cells.Item(76,4).Value = 1746800
$ws.Cells.Item(76,5).Value = 1666500
$ws.Cells.Item(77,4).Value = 0
$ws.Cells.Item(77,5).Value = 0
$ws.Cells.Item(80,4).Value = 43465
$ws.Cells.Item(80,5).Value = 43373
$ws.Cells.Item(81,4).Value = 91600
$ws.Cells.Item(81,5).Value = 57600
$ws.Cells.Item(83,4).Value = 125500
$ws.Cells.Item(83,5).Value = 107300
$ws.Cells.Item(84,4).Value = 0
$ws.Cells.Item(84,5).Value = 0
$ws.Cells.Item(85,4).Value = 0
$ws.Cells.Item(85,5).Value = 0
$ws.Cells.Item(86,4).Value = 0
$ws.Cells.Item(86,5).Value = 0
$ws.Cells.Item(87,4).Value = 0
$ws.Cells.Item(87,5).Value = 0
$ws.Cells.Item(88,4).Value = 0
$ws.Cells.Item(88,5).Value = 0
$ws.Cells.Item(89,4).Value = 216600
$ws.Cells.Item(89,5).Value = 181300
$ws.Cells.Item(91,4).Value = -4000
$ws.Cells.Item(91,5).Value = -9200
$ws.Cells.Item(92,4).Value = 0
$ws.Cells.Item(92,5).Value = 0
$ws.Cells.Item(93,4).Value = 0
$ws.Cells.Item(93,5).Value = 0
$ws.Cells.Item(94,4).Value = -219200
$ws.Cells.Item(94,5).Value = -212300
$ws.Cells.Item(96,4).Value = 0
$ws.Cells.Item(96,5).Value = 0
$ws.Cells.Item(97,4).Value = 0
$ws.Cells.Item(97,5).Value = 0
$ws.Cells.Item(98,4).Value = 0
$ws.Cells.Item(98,5).Value = 0
$ws.Cells.Item(99,4).Value = 0
$ws.Cells.Item(99,5).Value = 0
$ws.Cells.Item(100,4).Value = -36500
$ws.Cells.Item(100,5).Value = 236800
$ws.Cells.Item(101,4).Value = 0
$ws.Cells.Item(101,5).Value = 0
$ws.Cells.Item(102,4).Value = -39100
$ws.Cells.Item(102,5).Value = 205800

# Step 4: row 91 ("Capital Expenditures") was also restated for the five
#         quarters now sitting in F:J (the previously-reported values that
#         merely shifted right were replaced with corrected figures).
$ws.Cells.Item(91,6).Value = 4400
$ws.Cells.Item(91,7).Value = -7200
$ws.Cells.Item(91,8).Value = -15300
$ws.Cells.Item(91,9).Value = 3900
$ws.Cells.Item(91,10).Value = -2600
